$wb = $excel.ActiveWorkbook

# ALC sheet
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 2515
$ws.Range("J26").Value = 2515
$ws.Range("L26").Value = 2515
$ws.Range("N26").Value = -3203

$ws.Range("H76").Value = 4583.1665
$ws.Range("I76").Value = 4859.4
$ws.Range("K76").Value = 4859.4
$ws.Range("M76").Value = -4544.4

$ws.Range("H79").Value = 4583.1665
$ws.Range("I79").Value = 4859.4
$ws.Range("K79").Value = 4859.4
$ws.Range("M79").Value = -3767.4

# ARM sheet
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2086380.2
$ws.Range("I63").Value = 4763730.5
$ws.Range("J63").Value = 3996.7778
$ws.Range("K63").Value = 4763730.5
$ws.Range("L63").Value = 3996.7778
$ws.Range("M63").Value = -4763044.5
$ws.Range("N63").Value = -5368.7778

$ws.Range("H66").Value = 2086380.2
$ws.Range("I66").Value = 4763730.5
$ws.Range("J66").Value = 3996.7778
$ws.Range("K66").Value = 23818652.5
$ws.Range("L66").Value = 19983.889
$ws.Range("M66").Value = -23815220.5
$ws.Range("N66").Value = -26847.889

$ws.Range("H88").Value = 11234.267
$ws.Range("I88").Value = 1416.6666
$ws.Range("J88").Value = 17779.334
$ws.Range("K88").Value = 1416.6666
$ws.Range("L88").Value = 17779.334
$ws.Range("M88").Value = -1010.6666
$ws.Range("N88").Value = -18591.334

$ws.Range("H91").Value = 11234.267
$ws.Range("I91").Value = 1416.6666
$ws.Range("J91").Value = 17779.334
$ws.Range("K91").Value = 1416.6666
$ws.Range("L91").Value = 17779.334
$ws.Range("M91").Value = -12.66660000000002
$ws.Range("N91").Value = -20587.334

$ws.Range("H97").Value = 1460
$ws.Range("I97").Value = 1750
$ws.Range("J97").Value = 1266.6666
$ws.Range("K97").Value = 1750
$ws.Range("L97").Value = 1266.6666
$ws.Range("M97").Value = -1254
$ws.Range("N97").Value = -2258.6666

# BSM sheet
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3572759.5
$ws.Range("I86").Value = 1396
$ws.Range("J86").Value = 7693563.5
$ws.Range("K86").Value = 1396
$ws.Range("L86").Value = 7693563.5
$ws.Range("M86").Value = -273
$ws.Range("N86").Value = -7695809.5

$ws.Range("H89").Value = 3572759.5
$ws.Range("I89").Value = 1396
$ws.Range("J89").Value = 7693563.5
$ws.Range("K89").Value = 6980
$ws.Range("L89").Value = 38467817.5
$ws.Range("M89").Value = -1364
$ws.Range("N89").Value = -38479049.5

$ws.Range("H94").Value = 623.5
$ws.Range("I94").Value = 739.1667
$ws.Range("J94").Value = 450
$ws.Range("K94").Value = 739.1667
$ws.Range("L94").Value = 450
$ws.Range("M94").Value = -288.1667
$ws.Range("N94").Value = -1352

$ws.Range("H105").Value = 1731.8667
$ws.Range("I105").Value = 1759.8462
$ws.Range("J105").Value = 1550
$ws.Range("K105").Value = 1759.8462
$ws.Range("L105").Value = 1550
$ws.Range("M105").Value = -12.84619999999995
$ws.Range("N105").Value = -5044

# CRP sheet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5473808
$ws.Range("I31").Value = 6947524.5
$ws.Range("J31").Value = 32393.385
$ws.Range("K31").Value = 6947524.5
$ws.Range("L31").Value = 32393.385
$ws.Range("M31").Value = -6947229.5
$ws.Range("N31").Value = -32983.38499999999

$ws.Range("H34").Value = 5473808
$ws.Range("I34").Value = 6947524.5
$ws.Range("J34").Value = 32393.385
$ws.Range("K34").Value = 6947524.5
$ws.Range("L34").Value = 32393.385
$ws.Range("M34").Value = -6947322.5
$ws.Range("N34").Value = -32797.38499999999

$ws.Range("H62").Value = 2971.25
$ws.Range("I62").Value = 2385
$ws.Range("J62").Value = 3166.6667
$ws.Range("K62").Value = 2385
$ws.Range("L62").Value = 3166.6667
$ws.Range("M62").Value = -1761
$ws.Range("N62").Value = -4414.6667

$ws.Range("H65").Value = 2971.25
$ws.Range("I65").Value = 2385
$ws.Range("J65").Value = 3166.6667
$ws.Range("K65").Value = 11925
$ws.Range("L65").Value = 15833.3335
$ws.Range("M65").Value = -8805
$ws.Range("N65").Value = -22073.3335

$ws.Range("H134").Value = 4883748
$ws.Range("I134").Value = 4808499.5
$ws.Range("J134").Value = 5209823
$ws.Range("K134").Value = 14425498.5
$ws.Range("L134").Value = 15629469
$ws.Range("M134").Value = -14422963.5
$ws.Range("N134").Value = -15634539

# CUL sheet
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1122.3143
$ws.Range("I140").Value = 759.24
$ws.Range("J140").Value = 2030
$ws.Range("K140").Value = 2277.72
$ws.Range("L140").Value = 6090
$ws.Range("M140").Value = 2902.28
$ws.Range("N140").Value = -16450

# GSM sheet
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2079.3333
$ws.Range("I80").Value = 1762.8572
$ws.Range("J80").Value = 2356.25
$ws.Range("K80").Value = 1762.8572
$ws.Range("L80").Value = 2356.25
$ws.Range("M80").Value = -764.8571999999999
$ws.Range("N80").Value = -4352.25

$ws.Range("H83").Value = 2079.3333
$ws.Range("I83").Value = 1762.8572
$ws.Range("J83").Value = 2356.25
$ws.Range("K83").Value = 8814.286
$ws.Range("L83").Value = 11781.25
$ws.Range("M83").Value = -3822.286
$ws.Range("N83").Value = -21765.25

# LTW sheet
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2412.875
$ws.Range("I68").Value = 2835.3
$ws.Range("J68").Value = 1708.8334
$ws.Range("K68").Value = 2835.3
$ws.Range("L68").Value = 1708.8334
$ws.Range("M68").Value = -2086.3
$ws.Range("N68").Value = -3206.8334

$ws.Range("H71").Value = 2412.875
$ws.Range("I71").Value = 2835.3
$ws.Range("J71").Value = 1708.8334
$ws.Range("K71").Value = 14176.5
$ws.Range("L71").Value = 8544.166999999999
$ws.Range("M71").Value = -10432.5
$ws.Range("N71").Value = -16032.167

$ws.Range("H82").Value = 1935.069
$ws.Range("I82").Value = 1430.7778
$ws.Range("J82").Value = 2760.2727
$ws.Range("K82").Value = 1430.7778
$ws.Range("L82").Value = 2760.2727
$ws.Range("M82").Value = -1069.7778
$ws.Range("N82").Value = -3482.2727

$ws.Range("H85").Value = 1935.069
$ws.Range("I85").Value = 1430.7778
$ws.Range("J85").Value = 2760.2727
$ws.Range("K85").Value = 1430.7778
$ws.Range("L85").Value = 2760.2727
$ws.Range("M85").Value = -182.7778000000001
$ws.Range("N85").Value = -5256.2727

$ws.Range("H122").Value = 138893170
$ws.Range("I122").Value = 166669760
$ws.Range("K122").Value = 500009280
$ws.Range("M122").Value = -500006830

# WVR sheet
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54984
